$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("size_correction")
$ws2 = $wb.Worksheets.Item("scale_normalization")

# --- Update the Python-calculated isotopocule outputs (gamma/kappa inputs and
# --- the MATLAB/PYTHON OUTPUT block) on size_correction rows 3-17.
$ws1.Range("AP3").Value = 0.18304360488468446
$ws1.Range("AQ3").Value = 0.0915231290777685
$ws1.Range("AX3").Value = 14.9565698811999
$ws1.Range("AY3").Value = -1.30914416602279
$ws1.Range("AZ3").Value = 16.2657140472227
$ws1.Range("BA3").Value = 6.82371285758859
$ws1.Range("BB3").Value = 11.8952201489941
$ws1.Range("BC3").Value = 23.1813266136633
$ws1.Range("AP4").Value = 0.18304360488468446
$ws1.Range("AQ4").Value = 0.0915231290777685
$ws1.Range("AX4").Value = 6.73285280903357
$ws1.Range("AY4").Value = -12.4964138972027
$ws1.Range("AZ4").Value = 19.2292667062363
$ws1.Range("BA4").Value = -2.88178054408461
$ws1.Range("BB4").Value = 5.8070065868232
$ws1.Range("BC4").Value = 11.2845343634417
$ws1.Range("AP5").Value = 0.18304360488468446
$ws1.Range("AQ5").Value = 0.0915231290777685
$ws1.Range("AX5").Value = -0.692698082928267
$ws1.Range("AY5").Value = 2.68543047344715
$ws1.Range("AZ5").Value = -3.37812855637542
$ws1.Range("BA5").Value = 0.996366195259446
$ws1.Range("BB5").Value = 10.3652371768259
$ws1.Range("BC5").Value = 20.1852985766362
$ws1.Range("AP6").Value = 0.18304360488468446
$ws1.Range("AQ6").Value = 0.0915231290777685
$ws1.Range("AX6").Value = 14.6995121444024
$ws1.Range("AY6").Value = 0.545596282873539
$ws1.Range("AZ6").Value = 14.1539158615289
$ws1.Range("BA6").Value = 7.622554213638
$ws1.Range("BB6").Value = 11.9799369101096
$ws1.Range("BC6").Value = 23.3473440697933
$ws1.Range("AP7").Value = 0.18304360488468446
$ws1.Range("AQ7").Value = 0.0915231290777685
$ws1.Range("AX7").Value = 6.51596086226069
$ws1.Range("AY7").Value = -12.4175729116019
$ws1.Range("AZ7").Value = 18.9335337738626
$ws1.Range("BA7").Value = -2.95080602467062
$ws1.Range("BB7").Value = 5.89907464030892
$ws1.Range("BC7").Value = 11.4639402184801
$ws1.Range("AP8").Value = 0.18304360488468446
$ws1.Range("AQ8").Value = 0.0915231290777685
$ws1.Range("AX8").Value = -1.30772536726342
$ws1.Range("AY8").Value = 0.695789086729181
$ws1.Range("AZ8").Value = -2.0035144539926
$ws1.Range("BA8").Value = -0.305968140267121
$ws1.Range("BB8").Value = 9.43082474023748
$ws1.Range("BC8").Value = 18.3576143656101
$ws1.Range("AP9").Value = 0.18304360488468446
$ws1.Range("AQ9").Value = 0.0915231290777685
$ws1.Range("AX9").Value = 15.8877099673331
$ws1.Range("AY9").Value = -3.46601917305211
$ws1.Range("AZ9").Value = 19.3537291403852
$ws1.Range("BA9").Value = 6.2108453971405
$ws1.Range("BB9").Value = 11.013352053047
$ws1.Range("BC9").Value = 21.4539245146896
$ws1.Range("AP10").Value = 0.18304360488468446
$ws1.Range("AQ10").Value = 0.0915231290777685
$ws1.Range("AX10").Value = 14.352288139259
$ws1.Range("AY10").Value = -33.4592046076752
$ws1.Range("AZ10").Value = 47.8114927469343
$ws1.Range("BA10").Value = -9.55345823420811
$ws1.Range("BB10").Value = 10.7526616842892
$ws1.Range("BC10").Value = 20.9435550857681
$ws1.Range("AP11").Value = 0.18304360488468446
$ws1.Range("AQ11").Value = 0.0915231290777685
$ws1.Range("AX11").Value = 17.0031811376099
$ws1.Range("AY11").Value = -3.17532313527801
$ws1.Range("AZ11").Value = 20.1785042728879
$ws1.Range("BA11").Value = 6.91392900116594
$ws1.Range("BB11").Value = 10.8861024643158
$ws1.Range("BC11").Value = 21.2047848135454
$ws1.Range("AP12").Value = 0.18304360488468446
$ws1.Range("AQ12").Value = 0.0915231290777685
$ws1.Range("AX12").Value = 6.8691236337024
$ws1.Range("AY12").Value = -13.1908425064558
$ws1.Range("AZ12").Value = 20.0599661401582
$ws1.Range("BA12").Value = -3.16085943637672
$ws1.Range("BB12").Value = 5.21983650599189
$ws1.Range("BC12").Value = 10.1407240638105
$ws1.Range("AP13").Value = 0.18304360488468446
$ws1.Range("AQ13").Value = 0.0915231290777685
$ws1.Range("AX13").Value = -0.0696361091557795
$ws1.Range("AY13").Value = 1.62265568257535
$ws1.Range("AZ13").Value = -1.69229179173113
$ws1.Range("BA13").Value = 0.776509786709789
$ws1.Range("BB13").Value = 10.2313866904628
$ws1.Range("BC13").Value = 19.9233934898952
$ws1.Range("AP14").Value = 0.18304360488468446
$ws1.Range("AQ14").Value = 0.0915231290777685
$ws1.Range("AX14").Value = 14.2431406435554
$ws1.Range("AY14").Value = -2.27977397502809
$ws1.Range("AZ14").Value = 16.5229146185835
$ws1.Range("BA14").Value = 5.98168333426368
$ws1.Range("BB14").Value = 10.8426664922129
$ws1.Range("BC14").Value = 21.1197490321015
$ws1.Range("AP15").Value = 0.18304360488468446
$ws1.Range("AQ15").Value = 0.0915231290777685
$ws1.Range("AX15").Value = 16.5325251060171
$ws1.Range("AY15").Value = -2.9351313709296
$ws1.Range("AZ15").Value = 19.4676564769467
$ws1.Range("BA15").Value = 6.79869686754375
$ws1.Range("BB15").Value = 11.3793258579801
$ws1.Range("BC15").Value = 22.1706220531208
$ws1.Range("AP16").Value = 0.18304360488468446
$ws1.Range("AQ16").Value = 0.0915231290777685
$ws1.Range("AX16").Value = 6.96341945640721
$ws1.Range("AY16").Value = -13.5093226554305
$ws1.Range("AZ16").Value = 20.4727421118378
$ws1.Range("BA16").Value = -3.27295159951168
$ws1.Range("BB16").Value = 5.69517564432398
$ws1.Range("BC16").Value = 11.0666388211797
$ws1.Range("AP17").Value = 0.18304360488468446
$ws1.Range("AQ17").Value = 0.0915231290777685
$ws1.Range("AX17").Value = 13.6101546598585
$ws1.Range("AY17").Value = -26.9702125183365
$ws1.Range("AZ17").Value = 40.5803671781951
$ws1.Range("BA17").Value = -6.68002892923896
$ws1.Range("BB17").Value = 10.4675483235519
$ws1.Range("BC17").Value = 20.3855126219725

# --- scale_normalization: row 17 (S2_ref_2) formulas get rewritten as part of
# --- finishing the isotope-standards adjustable-parameter wiring.
$ws2.Range("B17").Formula = "=M4"
$ws2.Range("C17").Formula = "=N4"
$ws2.Range("F17").Formula = "=LN(D17)"
$ws2.Range("G17").Formula = "=LN(E17)"
$ws2.Range("H17").Formula = "=LN(B17)"
$ws2.Range("I17").Formula = "=LN(C17)"

# --- Restore selections: scale_normalization cursor moves to K22, then
# --- size_correction (the active/tabSelected sheet) cursor moves to AY28.
$ws2.Activate()
$ws2.Range("K22").Select()

$ws1.Activate()
$ws1.Range("AY28").Select()
